$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("B:C").Insert()
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"
$ws.Range("B2:C27").Value = "UN"
$ws.Range("B1:C27").NumberFormat = "General"
